$d = $word.ActiveDocument

# The new TODO bullet item goes right after the "TODO" heading, i.e.
# right before the first existing bullet ("Implement 'buddy space'").
# Inserting before that paragraph makes the new (empty) paragraph
# inherit its list/paragraph/run formatting (numPr, indentation,
# spacing, fonts, etc.) exactly, instead of inheriting the heading's.
$refPara = $d.Paragraphs(2)
$refPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "Implement getting build datetimestamp for composing build versio number automatically."
